$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Next_Cosmetic")

# H4: 507 -> 1027
$ws.Range("H4").Value = 1027

# E6: "Nyakuza Mask`nA hat in time" -> "Count Olaf`nA series of Unfortunate Events"
$ws.Range("E6").Value = "Count Olaf`nA series of Unfortunate Events"

# H6: "" -> image URL
$ws.Range("H6").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1104629198011768892/count_olaf_1.png"

# H8: image URL -> ""
$ws.Range("H8").Value = ""

# H9: "" -> image URL
$ws.Range("H9").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1104629198229884958/count_olaf_jacket_2.png"

# D12: "snekiecr8" -> "clown_noes666"
$ws.Range("D12").Value = "clown_noes666"

# E12: 28059068 -> 469660616
$ws.Range("E12").Value = 469660616

# F12: "616 days" -> "107 days"
$ws.Range("F12").Value = "107 days"

# D14: "snekie" -> "clown noes666"
$ws.Range("D14").Value = "clown noes666"

# D15: "Snekiecr8" -> "Clown_noes666"
$ws.Range("D15").Value = "Clown_noes666"
